$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.2642
$ws.Range("A3").Value = -22.2179
$ws.Range("E3").Value = 16.0312
$ws.Range("E12").Value = 17.39160000000001
$ws.Range("A14").Value = -21.719
$ws.Range("A16").Value = -21.75429999999999
$ws.Range("D18").Value = -8.8986
$ws.Range("A21").Value = -20.22949999999999
$ws.Range("A23").Value = -20.35609999999998
$ws.Range("D24").Value = -7.327700000000005
$ws.Range("E24").Value = 16.97950000000001
$ws.Range("A25").Value = -21.7509
$ws.Range("D25").Value = -7.432100000000003
$ws.Range("E25").Value = 16.70600000000001
$ws.Range("A26").Value = -21.31709999999997
$ws.Range("D27").Value = -8.922600000000003
$ws.Range("A29").Value = -20.97439999999998
$ws.Range("D30").Value = -7.670000000000001
$ws.Range("D31").Value = -8.819300000000002
$ws.Range("D39").Value = -8.366999999999997
$ws.Range("A40").Value = -20.2839
$ws.Range("E41").Value = 16.23979999999999
$ws.Range("D42").Value = -8.752699999999999
$ws.Range("D48").Value = -7.141699999999999
$ws.Range("E50").Value = 16.3535
$ws.Range("D51").Value = -7.498900000000001
$ws.Range("D52").Value = -7.512200000000003
$ws.Range("A53").Value = -21.93589999999999
$ws.Range("E53").Value = 16.96410000000002
$ws.Range("D55").Value = -9.032600000000004
$ws.Range("D56").Value = -7.8538
$ws.Range("E56").Value = 16.63490000000001
$ws.Range("A57").Value = -23.08210000000001
$ws.Range("D57").Value = -9.179800000000002
$ws.Range("E57").Value = 16.6657
$ws.Range("E58").Value = 16.27640000000002
$ws.Range("A59").Value = -22.5903
$ws.Range("D60").Value = -7.910599999999996
$ws.Range("E61").Value = 16.52430000000001
$ws.Range("E63").Value = 17.39730000000002
$ws.Range("E64").Value = 17.3619
$ws.Range("A65").Value = -21.81909999999999
$ws.Range("A69").Value = -21.63269999999998
$ws.Range("E70").Value = 17.33260000000001
$ws.Range("E72").Value = 16.99770000000001
$ws.Range("D73").Value = -7.574400000000002
$ws.Range("D74").Value = -8.584300000000002
$ws.Range("A79").Value = -20.4642
$ws.Range("A83").Value = -21.72969999999999
$ws.Range("E86").Value = 16.49710000000001
$ws.Range("D89").Value = -7.180199999999992
$ws.Range("E89").Value = 17.42620000000002
$ws.Range("D90").Value = -7.780300000000003
$ws.Range("A91").Value = -21.41760000000001
$ws.Range("D92").Value = -5.824000000000002
$ws.Range("A93").Value = -21.07969999999997
$ws.Range("E98").Value = 15.37820000000001
$ws.Range("A100").Value = -21.8933
$ws.Range("E100").Value = 16.56430000000001
$ws.Range("E102").Value = 16.19669999999999
